$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R20" (row 10) - update the "From" threshold in column C from 18 to 1
$ws.Range("C10").Value = 1
